# Apply the "tokenHolderRewards" feature update to Sheet1 of
# view_functions_grise.xlsx:
#   - E16's "get stake Details" reward-field list drops the 3 reward
#     columns and gains scrapeDay/stakeType.
#   - Row 17 is repurposed from "get reward scrape Day" /
#     getRewardScarpeDay / scrapeDay into "get Stake Reward Details" /
#     checkStakeRewards / the 3 reward fields (and grows to match the
#     Consolas code-style used by row 16's function-name cell).
#   - Rows 18/19's "day" countdown becomes "dayLeft".
#   - Three new rows (51-53) describing the new tokenHolderRewards
#     section are appended below the existing "graph statistics" block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- E16: reward/field list for "get stake Details" -----------------
$ws.Cells.Item(16, 5).Value = "startDays`nlockDays`nfinalDay`ncloseDay`nscrapeDay`nstakeType`nstakedAmount`npenaltyAmount`nisActive`nisMature"
$ws.Rows.Item(16).RowHeight = 150

# --- Row 17: repurposed into "get Stake Reward Details" -------------
$ws.Rows.Item(17).RowHeight = 45
$ws.Cells.Item(17, 2).Value = "get Stake Reward Details"
$ws.Cells.Item(17, 3).Value = "checkStakeRewards"
$ws.Cells.Item(17, 5).Value = "transcRewardAmount`nPenaltyRewardAmount`nreservoirRewardAmount"

# C17 switches from the plain "getShortTermSlotLeft"-style format to the
# Consolas / code look already used by C16 (checkStakeByID).
$ws.Cells.Item(16, 3).Copy() | Out-Null
$ws.Cells.Item(17, 3).PasteSpecial(-4122) | Out-Null

# E17 switches to the plain wrap-text format already used by D17/B17.
$ws.Cells.Item(17, 2).Copy() | Out-Null
$ws.Cells.Item(17, 5).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Rows 18/19: countdown label "day" -> "dayLeft" ------------------
$ws.Cells.Item(18, 5).Value = "dayLeft"
$ws.Cells.Item(19, 5).Value = "dayLeft"

# --- New rows 51-53: tokenHolderRewards section ----------------------
$ws.Cells.Item(51, 1).Value = "tokenHolderRewards"
$ws.Cells.Item(51, 2).Value = "view token Holder reward Amount"
$ws.Cells.Item(51, 3).Value = "viewTokenHolderTranscReward"
$ws.Cells.Item(51, 4).Value = "nothing"
$ws.Cells.Item(51, 5).Value = "rewardAmount"
$ws.Rows.Item(51).RowHeight = 30
# A51 needs the bold "section title" look used by A9/A28/A32/A35/A39.
$ws.Cells.Item(9, 1).Copy() | Out-Null
$ws.Cells.Item(51, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(52, 2).Value = "claim token holder reward on every 7th day from launch"
$ws.Cells.Item(52, 3).Value = "claimTokenHolderTranscReward"
$ws.Cells.Item(52, 4).Value = "nothing"
$ws.Cells.Item(52, 5).Value = "rewardAmount"
$ws.Rows.Item(52).RowHeight = 30

$ws.Cells.Item(53, 2).Value = "timer to check reward is ready for claim"
$ws.Cells.Item(53, 3).Value = "timeToClaimWeeklyReward"
$ws.Cells.Item(53, 4).Value = "nothing"
$ws.Cells.Item(53, 5).Value = "daysLeft(if zero them claimReward button will be enabled)"
$ws.Rows.Item(53).RowHeight = 30

# --- View state: scroll so the new rows are visible, selection on E56
$excel.ActiveWindow.ScrollRow = 38
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E56").Select()

Write-Host "tokenHolderRewards section applied"
